$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values (Régression Linéaire stays the same label)
$ws.Range("B2").Value = 5.12
$ws.Range("C2").Value = 2.56

# Insert a new row at position 3 so the old row 3 (Random Forest) becomes row 4
$ws.Rows.Item(3).Insert()

# New row 3: Clustering
$ws.Range("A3").Value = "Clustering"
$ws.Range("B3").Value = 3.29
$ws.Range("C3").Value = 1.65

# Row 4 keeps "Random Forest" label but with updated values
$ws.Range("A4").Value = "Random Forest"
$ws.Range("B4").Value = 0.06
$ws.Range("C4").Value = 0.03

# Row 5: another Random Forest entry
$ws.Range("A5").Value = "Random Forest"
$ws.Range("B5").Value = 0.08
$ws.Range("C5").Value = 0.04

# Row 6: Cross-Validation - Decision Tree
$ws.Range("A6").Value = "Cross-Validation - Decision Tree"
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
